# Applies the diff:
#  - slide 1 ("Date: " + "2012-04-06"): merge the two runs into a single
#    run whose text is "Date: 2012-04-06" (same visible text, but the
#    XML keeps only one <a:r> run instead of two).
#  - slide 3: move an 11-shape cluster (3 rectangles + their connectors,
#    plus 2 more rectangles/connectors lower down) by a uniform delta
#    (dx = -1296144 EMU, dy = +432048 EMU), leaving the title / legend
#    shapes untouched.
#
# Notes on precision:
#  Shape.Left / Shape.Top are expressed in points (1 pt = 12700 EMU) by
#  the PowerPoint object model. This host stores the point value as a
#  single-precision float and then floors (not rounds) the EMU
#  conversion, so naive "EMU/12700" literals can land 1 EMU short. The
#  literal point constants below were solved so that this round trip
#  reproduces the exact target EMU offsets from the diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: merge "Date: " and "2012-04-06" runs into a single run/text.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item("Untertitel 2")
$fullRange = $subtitle.TextFrame.TextRange
$fullText = $fullRange.Text
$marker = "Date: 2012-04-06"
$startPos = $fullText.IndexOf($marker) + 1
$dateRange = $fullRange.Characters($startPos, $marker.Length)
$dateRange.Text = $marker

# ---------------------------------------------------------------------
# Slide 3: reposition the 11 moved shapes.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Id -> (new Left pt, new Top pt), tuned so the stored EMU offsets match
# the diff exactly (see precision note above).
$moves = @{
    90  = @(36.8145,   224.6407)
    94  = @(269.2813,  224.6407)
    95  = @(269.2813,  292.6797)
    98  = @(345.82522, 252.9903)
    101 = @(48.1544,   292.6797)
    104 = @(93.5137,   252.9903)
    144 = @(150.21292, 292.6797)
    149 = @(93.5137,   252.9903)
    152 = @(204.0772,  252.9903)
    16  = @(138.8731,  355.0489)
    17  = @(204.0772,  321.0293)
}

for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shp = $s3.Shapes.Item($i)
    $id = $shp.Id
    if ($moves.ContainsKey($id)) {
        $xy = $moves[$id]
        $shp.Left = $xy[0]
        $shp.Top = $xy[1]
    }
}
